# Fix KW issues for Jahreswechsel in corona report
# Updates the "Impffortschritt" sheet's weekly vaccination figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Impffortschritt")

# Row 3: Gesamt
# B3/C3 look like plain numbers, so force them to be stored as text
# (matching the source workbook, where they are shared strings) and then
# drop the temporary "Text" number format back to the default style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "4463067"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "5777788"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "29,5 %"

# Row 4: davon in Impfzentren und Betrieben
$ws.Range("B4").Value = "1301719 ( 29,2 %)"
$ws.Range("C4").Value = "1695805 ( 29,4 %)"
$ws.Range("D4").Value = "30,3 %"

# Row 5: davon in ärztl. Praxen
$ws.Range("B5").Value = "3161348 ( 70,8 %)"
$ws.Range("C5").Value = "4081983 ( 70,6 %)"
